# "Updated data from BAG"
#
# 1) covid19_cases_switzerland sheet: fill in the newly-reported cantonal
#    case numbers for row 14 (date 43908) and move the selection.
# 2) demographics sheet: add the new L (pop/100000), M (test-positivity %),
#    N (estimated tests = ROUND(L*M,0)) and O (canton label) columns, and
#    move the selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: covid19_cases_switzerland
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("covid19_cases_switzerland")

$ws1.Range("B14").Value = 69
$ws1.Range("C14").Value = 3
$ws1.Range("D14").Value = 10
$ws1.Range("G14").Value = 207
$ws1.Range("H14").Value = 67
$ws1.Range("J14").Value = 10
$ws1.Range("K14").Value = 95
$ws1.Range("M14").Value = 48
$ws1.Range("P14").Value = 16
$ws1.Range("S14").Value = 31
$ws1.Range("T14").Value = 29
$ws1.Range("V14").Value = 506
$ws1.Range("X14").Value = 721
$ws1.Range("Y14").Value = 129
$ws1.Range("Z14").Value = 17
$ws1.Range("AB14").Value = 3028

$ws1.Activate() | Out-Null
$ws1.Range("S26").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: demographics
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("demographics")

# Per-canton test-positivity rate (%) used to derive estimated test counts.
$mValues = @{
    2  = 21.8
    3  = 15
    4  = 90.2
    5  = 10.2
    6  = 12.2
    7  = 50.7
    8  = 11.7
    9  = 143.2
    10 = 37.5
    11 = 21
    12 = 39.6
    13 = 6.9
    14 = 11.3
    15 = 47.9
    16 = 106.3
    17 = 43
    18 = 18.2
    19 = 13.4
    20 = 2.4
    21 = 30
    22 = 18.1
    23 = 34.7
    24 = 24.8
    25 = 42.3
    26 = 2.7
    27 = 18.6
}

# Canton label repeated in the new column O (mirrors column A).
$cantons = @{
    2  = "ZH"
    3  = "BE"
    4  = "VD"
    5  = "AG"
    6  = "SG"
    7  = "GE"
    8  = "LU"
    9  = "TI"
    10 = "VS"
    11 = "FR"
    12 = "BL"
    13 = "TG"
    14 = "SO"
    15 = "GR"
    16 = "BS"
    17 = "NE"
    18 = "SZ"
    19 = "ZG"
    20 = "SH"
    21 = "JU"
    22 = "AR"
    23 = "NW"
    24 = "GL"
    25 = "OW"
    26 = "UR"
    27 = "AI"
    28 = "CH"
}

# L2:L28 -> population / 100000 (single fill so Excel keeps it one shared formula).
$ws4.Range("L2:L28").Formula = "=B2/100000"

# M2:M27 -> literal test-positivity percentages entered per canton.
foreach ($r in 2..27) {
    $ws4.Cells.Item($r, 13).Value = $mValues[$r]
}

# N28 (=SUM) is set before the N2:N27 fill so the ROUND shared-formula block
# only ever covers N2:N27 and N28 keeps its own distinct SUM formula.
$ws4.Range("N28").Formula = "=SUM(N2:N27)"
$ws4.Range("N2:N27").Formula = "=ROUND(L2*M2,0)"

# O2:O28 -> canton label literal (same shared string as column A, not a formula).
foreach ($r in 2..28) {
    $ws4.Cells.Item($r, 15).Value = $cantons[$r]
}

$ws4.Activate() | Out-Null
$ws4.Range("N2:O27").Select() | Out-Null
$excel.ActiveWindow.RangeSelection.Item(1).Activate() | Out-Null
$ws4.Range("O27").Activate() | Out-Null
